# The paragraph containing "{m:userdoc 'zone1'}" currently holds two runs:
#   Run 1: "{m"
#   Run 2: ":userdoc 'zone1'}"
# The commit switches the parser to TokenIteratorFieldRewriterSplit, which
# rewrites the field token by emitting one run per logical token fragment:
#   "{" / "m" / ":userdoc 'zone1'" / "}"
# We locate the field text with Find, then rebuild that exact span as four
# separate <w:r> runs via WordOpenXML/InsertXML (the last run keeps
# xml:space="preserve" just like the migrated fixture).

$d = $word.ActiveDocument

$searchText = "{m:userdoc 'zone1'}"

$finder = $d.Range(0, 0)
$found = $finder.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $target = $d.Range($finder.Start, $finder.End)

    $openXml = '<?xml version="1.0" standalone="yes"?>' + `
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
        '<pkg:xmlData>' + `
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
        '<w:body>' + `
        '<w:p>' + `
        '<w:r><w:t>{</w:t></w:r>' + `
        '<w:r><w:t>m</w:t></w:r>' + `
        '<w:r><w:t>:userdoc ''zone1''</w:t></w:r>' + `
        '<w:r><w:t xml:space="preserve">}</w:t></w:r>' + `
        '</w:p>' + `
        '</w:body>' + `
        '</w:document>' + `
        '</pkg:xmlData>' + `
        '</pkg:part>' + `
        '</pkg:package>'

    $target.InsertXML($openXml)
}
